# Apply cell-value updates from the cryptos list refresh.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells receiving values that look numeric (e.g. "206.62", "0.990") must be
# pre-formatted as Text so Excel keeps them as literal strings (preserving
# exact formatting such as trailing zeros) instead of silently converting
# them to numbers.
$textCells = @("D5", "D17", "D19", "D23", "D25", "D26", "D27", "D30", "D31", "D35", "D37", "D39", "D42", "D44", "D46", "D48")
foreach ($addr in $textCells) {
    $ws.Range($addr).NumberFormat = "@"
}

# Row-by-row updates (Price column D, Volume(1h) column E).
$ws.Range("D2").Value = '26.869.94'
$ws.Range("E2").Value = '  -2.15%  '
$ws.Range("D3").Value = '1.567.30'
$ws.Range("E3").Value = '  -0.21%  '
$ws.Range("E4").Value = '  -0.04%  '
$ws.Range("D5").Value = '206.62'
$ws.Range("E5").Value = '  -1.11%  '
$ws.Range("E6").Value = '  -1.82%  '
$ws.Range("E7").Value = '  +0.00%  '
$ws.Range("E9").Value = '  -0.43%  '
$ws.Range("E10").Value = '  -1.09%  '
$ws.Range("E11").Value = '  -0.27%  '
$ws.Range("D12").Value = '1.790.28'
$ws.Range("E12").Value = '  -0.21%  '
$ws.Range("D13").Value = '1.559.39'
$ws.Range("E13").Value = '  -0.88%  '
$ws.Range("E14").Value = '  -2.13%  '
$ws.Range("E15").Value = '  -0.63%  '
$ws.Range("D16").Value = '26.860.83'
$ws.Range("E16").Value = '  -2.19%  '
$ws.Range("D17").Value = '61.53'
$ws.Range("E17").Value = '  -3.56%  '
$ws.Range("E18").Value = '  +1.62%  '
$ws.Range("D19").Value = '215.10'
$ws.Range("E19").Value = '  +0.61%  '
$ws.Range("E20").Value = '  -2.06%  '
$ws.Range("E21").Value = '  +0.01%  '
$ws.Range("E22").Value = '  -0.07%  '
$ws.Range("D23").Value = '9.33'
$ws.Range("E23").Value = '  -2.57%  '
$ws.Range("E24").Value = '  -0.96%  '
$ws.Range("D25").Value = '154.12'
$ws.Range("E25").Value = '  +1.10%  '
$ws.Range("D26").Value = '6.73'
$ws.Range("E26").Value = '  +0.03%  '
$ws.Range("D27").Value = '14.98'
$ws.Range("E29").Value = '  -1.11%  '
$ws.Range("D30").Value = '0.0467'
$ws.Range("E30").Value = '  -0.78%  '
$ws.Range("D31").Value = '1.12'
$ws.Range("E31").Value = '  -3.33%  '
$ws.Range("E32").Value = '  -0.76%  '
$ws.Range("D33").Value = '1.403.19'
$ws.Range("E33").Value = '  +1.50%  '
$ws.Range("D35").Value = '1.53'
$ws.Range("E35").Value = '  -1.14%  '
$ws.Range("E36").Value = '  -1.01%  '
$ws.Range("D37").Value = '0.935'
$ws.Range("E37").Value = '  -1.76%  '
$ws.Range("E38").Value = '  -2.52%  '
$ws.Range("D39").Value = '0.528'
$ws.Range("E39").Value = '  -2.87%  '
$ws.Range("E40").Value = '  -1.44%  '
$ws.Range("E41").Value = '  +0.01%  '
$ws.Range("D42").Value = '0.990'
$ws.Range("E42").Value = '  +0.66%  '
$ws.Range("E43").Value = '  -0.29%  '
$ws.Range("D44").Value = '5.33'
$ws.Range("E44").Value = '  +1.09%  '
$ws.Range("E45").Value = '  +0.58%  '
$ws.Range("D46").Value = '63.32'
$ws.Range("E46").Value = '  -1.50%  '
$ws.Range("D47").Value = '1.703.02'
$ws.Range("E47").Value = '  -0.05%  '
$ws.Range("D48").Value = '86.13'
$ws.Range("E48").Value = '  +0.64%  '
$ws.Range("D49").Value = '0.0₇0984'
$ws.Range("E49").Value = '  -1.53%  '
$ws.Range("E50").Value = '  -0.37%  '
$ws.Range("E51").Value = '  -0.93%  '

# Restore the default (General) formatting on the cells we forced to Text,
# without disturbing the text values that are now stored in them.
foreach ($addr in $textCells) {
    $ws.Range($addr).ClearFormats()
}

